$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2272727272727273
$ws.Range("C2").Value = 0.5064935064935064
$ws.Range("J2").Value = 0.01948051948051948
$ws.Range("O2").Value = 0.006493506493506494
$ws.Range("P2").Value = 0.1623376623376623
$ws.Range("S2").Value = 0.07792207792207792
$ws.Range("B3").Value = 0.01829268292682927
$ws.Range("C3").Value = 0.0426829268292683
$ws.Range("J3").Value = 0.04878048780487805
$ws.Range("P3").Value = 0.676829268292683
$ws.Range("S3").Value = 0.2134146341463415
$ws.Range("J4").Value = 0.06976744186046512
$ws.Range("P4").Value = 0.6511627906976745
$ws.Range("S4").Value = 0.2790697674418605
$ws.Range("B6").Value = 0.0371900826446281
$ws.Range("D6").Value = 0.01239669421487603
$ws.Range("F6").Value = 0.06198347107438017
$ws.Range("J6").Value = 0.2809917355371901
$ws.Range("O6").Value = 0.02479338842975207
$ws.Range("Q6").Value = 0.1074380165289256
$ws.Range("R6").Value = 0.08264462809917356
$ws.Range("S6").Value = 0.3925619834710744
$ws.Range("B7").Value = 0.1005025125628141
$ws.Range("D7").Value = 0.01005025125628141
$ws.Range("F7").Value = 0.08040201005025126
$ws.Range("J7").Value = 0.1256281407035176
$ws.Range("O7").Value = 0.02010050251256281
$ws.Range("Q7").Value = 0.1557788944723618
$ws.Range("R7").Value = 0.1155778894472362
$ws.Range("S7").Value = 0.3919597989949749
$ws.Range("B8").Value = 0.09108910891089109
$ws.Range("D8").Value = 0.01188118811881188
$ws.Range("F8").Value = 0.08514851485148515
$ws.Range("J8").Value = 0.1207920792079208
$ws.Range("O8").Value = 0.0297029702970297
$ws.Range("Q8").Value = 0.1326732673267327
$ws.Range("R8").Value = 0.08118811881188119
$ws.Range("S8").Value = 0.4475247524752475
$ws.Range("B9").Value = 0.0738255033557047
$ws.Range("D9").Value = 0.02684563758389262
$ws.Range("F9").Value = 0.0738255033557047
$ws.Range("J9").Value = 0.1342281879194631
$ws.Range("O9").Value = 0.02684563758389262
$ws.Range("Q9").Value = 0.1342281879194631
$ws.Range("R9").Value = 0.1208053691275168
$ws.Range("S9").Value = 0.4093959731543624
$ws.Range("B10").Value = 0.11198738170347
$ws.Range("D10").Value = 0.02287066246056782
$ws.Range("E10").Value = 0.0007886435331230284
$ws.Range("F10").Value = 0.0694006309148265
$ws.Range("J10").Value = 0.110410094637224
$ws.Range("O10").Value = 0.02129337539432177
$ws.Range("Q10").Value = 0.1916403785488959
$ws.Range("R10").Value = 0.07649842271293375
$ws.Range("S10").Value = 0.3951104100946372
$ws.Range("G11").Value = 0.1627906976744186
$ws.Range("J11").Value = 0.09302325581395349
$ws.Range("K11").Value = 0.2225913621262458
$ws.Range("L11").Value = 0.5083056478405316
$ws.Range("S11").Value = 0.0132890365448505
$ws.Range("G12").Value = 0.8089171974522293
$ws.Range("J12").Value = 0.1401273885350318
$ws.Range("L12").Value = 0.02547770700636943
$ws.Range("S12").Value = 0.02547770700636943
$ws.Range("G13").Value = 0.7804878048780488
$ws.Range("J13").Value = 0.1951219512195122
$ws.Range("S13").Value = 0.02439024390243903
$ws.Range("F15").Value = 0.00966183574879227
$ws.Range("H15").Value = 0.1545893719806763
$ws.Range("I15").Value = 0.03381642512077294
$ws.Range("J15").Value = 0.3478260869565217
$ws.Range("K15").Value = 0.0821256038647343
$ws.Range("M15").Value = 0.004830917874396135
$ws.Range("O15").Value = 0.03864734299516908
$ws.Range("S15").Value = 0.3285024154589372
$ws.Range("F16").Value = 0.01081081081081081
$ws.Range("H16").Value = 0.2486486486486487
$ws.Range("I16").Value = 0.03243243243243243
$ws.Range("J16").Value = 0.3837837837837838
$ws.Range("K16").Value = 0.1189189189189189
$ws.Range("M16").Value = 0.01621621621621622
$ws.Range("O16").Value = 0.05405405405405406
$ws.Range("S16").Value = 0.1351351351351351
$ws.Range("F17").Value = 0.02110817941952507
$ws.Range("H17").Value = 0.1794195250659631
$ws.Range("I17").Value = 0.079155672823219
$ws.Range("J17").Value = 0.4379947229551451
$ws.Range("K17").Value = 0.079155672823219
$ws.Range("M17").Value = 0.01846965699208443
$ws.Range("N17").Value = 0.002638522427440633
$ws.Range("O17").Value = 0.05277044854881267
$ws.Range("S17").Value = 0.129287598944591
$ws.Range("F18").Value = 0.02985074626865672
$ws.Range("H18").Value = 0.1940298507462687
$ws.Range("I18").Value = 0.05472636815920398
$ws.Range("J18").Value = 0.3781094527363184
$ws.Range("K18").Value = 0.1194029850746269
$ws.Range("M18").Value = 0.01492537313432836
$ws.Range("O18").Value = 0.05970149253731343
$ws.Range("S18").Value = 0.1492537313432836
$ws.Range("F19").Value = 0.01531728665207877
$ws.Range("H19").Value = 0.2319474835886214
$ws.Range("I19").Value = 0.0700218818380744
$ws.Range("J19").Value = 0.3778264040846098
$ws.Range("K19").Value = 0.1021152443471918
$ws.Range("M19").Value = 0.02042304886943837
$ws.Range("O19").Value = 0.05470459518599562
$ws.Range("S19").Value = 0.1276440554339898
